# Update "想去人数" (number of people interested) values in column F
# for rows 2-5 on the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

$values = @{
    2 = 42
    3 = 141
    4 = 15
    5 = 38
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $values.Keys) {
        $ws.Range("F$row").Value = $values[$row]
    }
}
